# Fix multiple bugs in the Meeting Part Class sheet ("Input List"):
#  - Split the single "Chairman" column into two columns: "Chairman Tue" and
#    "Chairman Fri" (the new "Fri" column duplicates the Tue assignments).
#  - This requires inserting a new column before column B, which shifts the
#    former B:I columns to C:J.
#  - Remove the stray highlight formatting that had been left on the
#    "Treasures 2" cells for the Fri row (old B6:C6, now C6:D6).
#  - Update the HTML_all / HTML_tables defined names so they keep pointing at
#    the "Enter the start date below" helper column after the column shift.
#  - Leave the active selection on B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column B; this shifts old B..I to C..J and
#    carries formatting from the left neighbour (column A) into new column B.
$ws.Columns("B:B").Insert()

# 2. Split the header: A1 becomes the "Tue" chairman column, new B1 becomes
#    the "Fri" chairman column.
$ws.Range("A1").Value = "Chairman Tue"
$ws.Range("B1").Value = "Chairman Fri"

# 3. The new "Fri" column starts out as a copy of the "Tue" column values.
$ws.Range("A2:A10").Copy()
$ws.Range("B2:B10").PasteSpecial(-4163)

# 4. Clear the stray highlighted-box formatting that had been sitting on the
#    "Treasures 2" Fri cells (now C6:D6 after the column insert) by copying
#    the normal format from A6.
$ws.Range("A6").Copy()
$ws.Range("C6:D6").PasteSpecial(-4122)

# 5. Fix up the HTML_all / HTML_tables defined names so they still refer to
#    the "Enter the start date below" cell, which moved from F1 to G1.
$wb.Names("HTML_all").RefersTo = "='Input List'!`$G`$1:`$G`$1"
$wb.Names("HTML_tables").RefersTo = "='Input List'!`$G`$1:`$G`$1"

# 6. Leave the selection on the newly added B2 cell.
[void]$ws.Range("B2").Select()
